$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.939.63'
$ws.Range('E2').Value = '  +2.86%  '
$ws.Range('D3').Value = '3.805.83'
$ws.Range('E3').Value = '  +1.07%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '698.98'
$ws.Range('E5').Value = '  +8.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.52'
$ws.Range('E6').Value = '  +4.22%  '
$ws.Range('D7').Value = '3.804.97'
$ws.Range('E7').Value = '  +1.05%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +1.28%  '
$ws.Range('E10').Value = '  +2.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.39'
$ws.Range('E11').Value = '  +7.18%  '
$ws.Range('E12').Value = '  +0.99%  '
$ws.Range('E13').Value = '  +7.96%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.45'
$ws.Range('D15').Value = '4.442.35'
$ws.Range('E15').Value = '  +0.96%  '
$ws.Range('D16').Value = '3.823.50'
$ws.Range('E16').Value = '  +1.74%  '
$ws.Range('D17').Value = '70.905.12'
$ws.Range('E17').Value = '  +2.87%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.86'
$ws.Range('E18').Value = '  +1.13%  '
$ws.Range('E19').Value = '  +2.97%  '
$ws.Range('E20').Value = '  +0.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.09'
$ws.Range('E21').Value = '  +16.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '482.70'
$ws.Range('E22').Value = '  +2.32%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.715'
$ws.Range('E23').Value = '  +1.65%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.31'
$ws.Range('E24').Value = '  +2.97%  '
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.44'
$ws.Range('E26').Value = '  +2.43%  '
$ws.Range('E27').Value = '  +3.72%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.51'
$ws.Range('E28').Value = '  +4.41%  '
$ws.Range('D29').Value = '3.955.69'
$ws.Range('E29').Value = '  +1.05%  '
$ws.Range('E31').Value = '  +15.83%  '
$ws.Range('E32').Value = '  +6.17%  '
$ws.Range('E33').Value = '  +1.44%  '
$ws.Range('E34').Value = '  +6.05%  '
$ws.Range('E35').Value = '  +3.52%  '
$ws.Range('E36').Value = '  +4.85%  '
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('E38').Value = '  +2.54%  '
$ws.Range('E39').Value = '  +7.08%  '
$ws.Range('E40').Value = '  +4.61%  '
$ws.Range('E41').Value = '  +12.45%  '
$ws.Range('B42').Value = 'Mantle'
$ws.Range('C42').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.978'
$ws.Range('E42').Value = '  +2.25%  '
$ws.Range('B43').Value = 'FLOKI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.000328'
$ws.Range('E43').Value = '  +23.52%  '
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '162.38'
$ws.Range('E46').Value = '  +4.59%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '44.97'
$ws.Range('E47').Value = '  +0.28%  '
$ws.Range('E48').Value = '  +2.70%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.302'
$ws.Range('E49').Value = '  +2.58%  '
$ws.Range('E50').Value = '  -1.84%  '
$ws.Range('E51').Value = '  +2.71%  '

Write-Host "Applied 75 cell updates"
